$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 6 (year 2025) metrics per latest data refresh
$ws.Range("C6").Value = 412
$ws.Range("E6").Value = 105
$ws.Range("G6").Value = 25.48543689320388
$ws.Range("H6").Value = 74.51456310679612
